$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1470
$ws1.Range("F8").Value = 6181
$ws1.Range("F12").Value = 5060
$ws1.Range("F16").Value = 54
$ws1.Range("F22").Value = 3579
$ws1.Range("F23").Value = 148

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1470
$ws4.Range("F9").Value = 6181
$ws4.Range("F13").Value = 5060
$ws4.Range("F17").Value = 54
$ws4.Range("F23").Value = 3579
$ws4.Range("F25").Value = 148
